$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6 through 21 (entire rows), leaving header + 4 data rows
$ws.Range("A6:B21").EntireRow.Delete() | Out-Null

# Update the remaining data rows (2-5) with the new id / name values
$ws.Range("A2").Value = 484931
$ws.Range("B2").Value = "Шустер В.Л., Шустер Владимир Львович"

$ws.Range("A3").Value = 534346
$ws.Range("B3").Value = "Пунанова С.А., Пунанова Светлана Александровна"

$ws.Range("A4").Value = 1090961
$ws.Range("B4").Value = "Колоколова И.В., Колоколова Ирина Владимировна"

$ws.Range("A5").Value = 9036402995
$ws.Range("B5").Value = "Краус З.Т., Краус Зоя Тимофеевна"
